$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 05:47"

# Row 29: Kazajistan - updated case counts
$ws.Range("B29").Value = 87664
$ws.Range("C29").Value = 1472
$ws.Range("D29").Value = 57815
$ws.Range("E29").Value = 29056

# Row 39: Belgica - updated case counts
$ws.Range("B39").Value = 67335
$ws.Range("C39").Value = 673
$ws.Range("D39").Value = 17491
$ws.Range("E39").Value = 40008
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 9836

# Row 51: Honduras - updated case counts
$ws.Range("B51").Value = 40944
$ws.Range("C51").Value = 484
$ws.Range("D51").Value = 5281
$ws.Range("E51").Value = 34404
$ws.Range("G51").Value = 45
$ws.Range("H51").Value = 1259

# Row 91: Haiti - updated case counts
$ws.Range("B91").Value = 7378
$ws.Range("C91").Value = 7
$ws.Range("E91").Value = 2752
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 159

# Rows 183/184: Monaco moved above Aruba in the country list (order swap),
# along with each country's updated case counts.
$ws.Range("A183").Value = "Monaco"
$ws.Range("B183").Value = 120
$ws.Range("D183").Value = 105
$ws.Range("E183").Value = 11
$ws.Range("H183").Value = 4

$ws.Range("A184").Value = "Aruba"
$ws.Range("B184").Value = 119
$ws.Range("D184").Value = 102
$ws.Range("E184").Value = 14
$ws.Range("H184").Value = 3
